$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Total" formula column C for rows 1 through 11 as a shared formula
$ws.Range("C1:C11").Formula = "=A1+B1"

# Update selection to match the new used range (C1:C11 instead of C1:C12)
$ws.Range("C1:C11").Select()
